$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new column header "DBLP Code" in D1, matching the bold header style of A1:C1
$ws.Range("D1").Value = "DBLP Code"
$ws.Range("D1").Font.Bold = $true

# Move the active selection to D1 (matches the post-edit cursor position)
$ws.Range("D1").Select()
